$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 90. Excel shifts rows 90..139 down to 91..140,
# preserving all of their existing values/formatting, and expands the used range.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly price record.
$ws.Cells.Item(90, 1).Value2 = 10
$ws.Cells.Item(90, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(90, 3).Value2 = "La Araucanía"
$ws.Cells.Item(90, 4).Value2 = 44830
$ws.Cells.Item(90, 5).Value2 = 9
$ws.Cells.Item(90, 6).Value2 = 100112035
$ws.Cells.Item(90, 7).Value2 = "Bruselas (repollito)"
$ws.Cells.Item(90, 8).Value2 = "Sin especificar"
$ws.Cells.Item(90, 9).Value2 = "Primera"
$ws.Cells.Item(90, 10).Value2 = 55
$ws.Cells.Item(90, 11).Value2 = 23000
$ws.Cells.Item(90, 12).Value2 = 23000
$ws.Cells.Item(90, 13).Value2 = 23000
$ws.Cells.Item(90, 14).Value2 = "$/malla 10 kilos"
$ws.Cells.Item(90, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(90, 16).Value2 = 2300
$ws.Cells.Item(90, 17).Value2 = 10
$ws.Cells.Item(90, 18).Value2 = "Hortaliza"
